# Generate Report for Handback
#
# This applies the "handback report" update to the localization-status
# workbook: the overview status text moves from "Ready for handoff" to
# "Handed back: in sync with en-US", the zh-cn/de-de sheets get their
# "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns filled in (with a hyperlink on the target-file cell,
# mirroring the source-file hyperlink), and several columns are widened
# so the newly-populated file-name / datetime columns are readable.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)
$zhcn     = $wb.Worksheets.Item(2)
$dede     = $wb.Worksheets.Item(3)

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (every cell currently showing the old status moves to the new one)
# ---------------------------------------------------------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. Overview sheet: widen the zh-cn / de-de status columns (E, F) so the
#    longer status text fits.
# ---------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 29.17
$overview.Columns.Item(6).ColumnWidth = 29.17

# ---------------------------------------------------------------------
# 3. zh-cn sheet (Worksheets.Item(2)): fill in handback report columns
# ---------------------------------------------------------------------
$zhcn.Columns.Item(3).ColumnWidth  = 29.17
$zhcn.Columns.Item(9).ColumnWidth  = 39.17
$zhcn.Columns.Item(10).ColumnWidth = 39.17

# Rebuild the hyperlinks in row order (A2, I2, A3, I3) so the new
# "Latest Target File" links land next to their source-file counterparts.
$zhcn.Hyperlinks.Delete()

$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1382505e488e8cceb86bc9659d0285d41e115685/e2e/06643d50-3e86-45a5-836e-544013d5253d.md", "", "", "06643d50-3e86-45a5-836e-544013d5253d.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1382505e488e8cceb86bc9659d0285d41e115685/e2e/06643d50-3e86-45a5-836e-544013d5253d.md", "", "", "06643d50-3e86-45a5-836e-544013d5253d.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1382505e488e8cceb86bc9659d0285d41e115685/e2e/6da22061-c3bb-4a6b-a302-ab5e396ff493.md", "", "", "6da22061-c3bb-4a6b-a302-ab5e396ff493.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1382505e488e8cceb86bc9659d0285d41e115685/e2e/6da22061-c3bb-4a6b-a302-ab5e396ff493.md", "", "", "6da22061-c3bb-4a6b-a302-ab5e396ff493.md") | Out-Null

# Latest Handback File
$zhcn.Range("J2").Value = "06643d50-3e86-45a5-836e-544013d5253d.3577508bacf7723954ec174b36cfa0a718ac7757.zh-cn.xlf"
$zhcn.Range("J3").Value = "6da22061-c3bb-4a6b-a302-ab5e396ff493.aecae32973175cb4369edc2f828e728be2c1b7b6.zh-cn.xlf"

# Latest Handback DateTime
$zhcn.Range("K2").Value = "2016-08-21 15:08:06"
$zhcn.Range("K3").Value = "2016-08-21 15:08:06"

# ---------------------------------------------------------------------
# 4. de-de sheet (Worksheets.Item(3)): fill in handback report columns
# ---------------------------------------------------------------------
$dede.Columns.Item(3).ColumnWidth  = 29.17
$dede.Columns.Item(9).ColumnWidth  = 39.17
$dede.Columns.Item(10).ColumnWidth = 39.17

$dede.Hyperlinks.Delete()

$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1382505e488e8cceb86bc9659d0285d41e115685/e2e/06643d50-3e86-45a5-836e-544013d5253d.md", "", "", "06643d50-3e86-45a5-836e-544013d5253d.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1382505e488e8cceb86bc9659d0285d41e115685/e2e/06643d50-3e86-45a5-836e-544013d5253d.md", "", "", "06643d50-3e86-45a5-836e-544013d5253d.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1382505e488e8cceb86bc9659d0285d41e115685/e2e/6da22061-c3bb-4a6b-a302-ab5e396ff493.md", "", "", "6da22061-c3bb-4a6b-a302-ab5e396ff493.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1382505e488e8cceb86bc9659d0285d41e115685/e2e/6da22061-c3bb-4a6b-a302-ab5e396ff493.md", "", "", "6da22061-c3bb-4a6b-a302-ab5e396ff493.md") | Out-Null

# Latest Handback File
$dede.Range("J2").Value = "06643d50-3e86-45a5-836e-544013d5253d.3577508bacf7723954ec174b36cfa0a718ac7757.de-de.xlf"
$dede.Range("J3").Value = "6da22061-c3bb-4a6b-a302-ab5e396ff493.aecae32973175cb4369edc2f828e728be2c1b7b6.de-de.xlf"

# Latest Handback DateTime
$dede.Range("K2").Value = "2016-08-21 15:08:13"
$dede.Range("K3").Value = "2016-08-21 15:08:13"

Write-Host "Handback report generated."
